$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 7 that is a copy of the current row 6 (old data),
# shifting nothing else. We'll set row 7 values directly to match the
# old row 6 contents, then overwrite row 6 with the new data.

# First, write the new row 7 (copy of the original row 6 data)
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = 44623
$ws.Range("D7").NumberFormat = $ws.Range("D6").NumberFormat
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112017
$ws.Range("G7").Value = "Corazón de apio"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 1800
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = 1900
$ws.Range("N7").Value = "$/paquete"
$ws.Range("O7").Value = "Región de Arica y Parinacota"
$ws.Range("P7").Value = 1900
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"

# Now overwrite row 6 with the new data
$ws.Range("D6").Value = 45218
$ws.Range("J6").Value = 180
$ws.Range("K6").Value = 1400
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = 1444
$ws.Range("N6").Value = "$/docena de matas"
$ws.Range("P6").Value = 241
$ws.Range("Q6").Value = 6
